$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, centered, bordered) used by the rest of row 1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill the season record for every data row (2-52) with the team's record
$ws.Range("AD2:AD52").Value = 73
$ws.Range("AE2:AE52").Value = 89
$ws.Range("AF2:AF52").Value = 0
